# HU Registrar paciente.docx edit
#
# 1) Merge the two runs "En caso de que" + "/Cuando" into a single run.
# 2) Swap the contents of the "Entonces" row and the following "Y" row:
#      - "Entonces" row gets the short text "El sistema valida que el
#        animal no se encuentre registrado" (split as "El " + rest, to
#        match the target run layout).
#      - "Y" row gets the long text that used to live in "Entonces"
#        (split into the same four runs it originally had).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# Change 1: merge "En caso de que" + "/Cuando" into one run
# ---------------------------------------------------------------------
$d.Content.Find.Execute("En caso de que/Cuando", $false, $false, $false, `
    $false, $false, $true, 1, $false, "En caso de que/Cuando", 2)

# ---------------------------------------------------------------------
# Change 2: "Entonces" row (row 11) -> short text, split into 2 runs
# ---------------------------------------------------------------------
$cellEntonces = $t.Cell(11, 2)
$rEntonces = $cellEntonces.Range
$rEntoncesBody = $d.Range($rEntonces.Start, $rEntonces.End - 1)
$rEntoncesBody.Text = "El sistema valida que el animal no se encuentre registrado"

# Force a run boundary after "El " by toggling a character property on
# and back off; the writer keeps the run split even once the formatting
# difference is removed again.
$cellEntonces = $t.Cell(11, 2)
$rEntonces = $cellEntonces.Range
$boundary = $d.Range($rEntonces.Start, $rEntonces.Start + 3)
$boundary.Bold = 1

$cellEntonces = $t.Cell(11, 2)
$rEntonces = $cellEntonces.Range
$boundary = $d.Range($rEntonces.Start, $rEntonces.Start + 3)
$boundary.Bold = 0

# ---------------------------------------------------------------------
# Change 3: "Y" row (row 12) -> long text, split into the original 4 runs
# ---------------------------------------------------------------------
$cellY = $t.Cell(12, 2)
$rY = $cellY.Range
$rYBody = $d.Range($rY.Start, $rY.End - 1)
$rYBody.Text = "El dueño del animal debe proporcionar la información requerida sobre su mascota para poder registrarla en el sistema (como su nombre, edad, raza, peso y cualquier afección médica existente) al igual que los datos del dueño (nombre, identificación, correo y dirección)"

# Run boundaries (character offsets from the start of the cell text):
#   0   .. 116 -> "El dueño ... en el sistema"
#   116 .. 118 -> " ("
#   118 .. 188 -> "como su nombre ... existente"
#   188 .. 267 -> ") al igual ... direccion)"
$offsets = @(116, 118, 188)
foreach ($off in $offsets) {
    $cellY = $t.Cell(12, 2)
    $rY = $cellY.Range
    $boundary = $d.Range($rY.Start, $rY.Start + $off)
    $boundary.Bold = 1

    $cellY = $t.Cell(12, 2)
    $rY = $cellY.Range
    $boundary = $d.Range($rY.Start, $rY.Start + $off)
    $boundary.Bold = 0
}

Write-Output "edit complete"
